$wb = $excel.ActiveWorkbook

# --- Sheet 1: Summary ---
$ws1 = $wb.Worksheets.Item("Summary")
$ws1.Range("B2").Value = 0.2117437722419929
$ws1.Range("C2").Value = 0.05376344086021505
$ws1.Range("D2").Value = 0.8928571428571429
$ws1.Range("E2").Value = 0.101419878296146
$ws1.Range("F2").Value = 0.2166377816291161
$ws1.Range("G2").Value = 0.5579399141630901
$ws1.Range("H2").Value = 0.79303772070626
$ws1.Range("I2").Value = 25
$ws1.Range("J2").Value = 440
$ws1.Range("K2").Value = 94
$ws1.Range("L2").Value = 3

# --- Sheet 2: Classification Report ---
$ws2 = $wb.Worksheets.Item("Classification Report")

$ws2.Range("B2").Value = 0.9690721649484536
$ws2.Range("C2").Value = 0.1760299625468165
$ws2.Range("D2").Value = 0.2979397781299525

$ws2.Range("B3").Value = 0.05376344086021505
$ws2.Range("C3").Value = 0.8928571428571429
$ws2.Range("D3").Value = 0.101419878296146

$ws2.Range("B4").Value = 0.2117437722419929
$ws2.Range("C4").Value = 0.2117437722419929
$ws2.Range("D4").Value = 0.2117437722419929
$ws2.Range("E4").Value = 0.2117437722419929

$ws2.Range("B5").Value = 0.5114178029043344
$ws2.Range("C5").Value = 0.5344435527019797
$ws2.Range("D5").Value = 0.1996798282130492

$ws2.Range("B6").Value = 0.9234695950650539
$ws2.Range("C6").Value = 0.2117437722419929
$ws2.Range("D6").Value = 0.2881487510919692

# --- Sheet 3: Confusion Matrix ---
$ws3 = $wb.Worksheets.Item("Confusion Matrix")

$ws3.Range("B2").Value = 94
$ws3.Range("C2").Value = 440

$ws3.Range("B3").Value = 3
$ws3.Range("C3").Value = 25
